$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45185,4497.97998046875,4497.97998046875,4447.2099609375,4450.31982421875,4450.31982421875,6932230000),
    @(45186,4497.97998046875,4497.97998046875,4447.2099609375,4450.31982421875,4450.31982421875,6932230000),
    @(45187,4445.1298828125,4466.35986328125,4442.10986328125,4453.52978515625,4453.52978515625,3161230000),
    @(45188,4445.41015625,4449.85009765625,4416.60986328125,4443.9501953125,4443.9501953125,3614880000),
    @(45189,4452.81005859375,4461.02978515625,4401.3798828125,4402.2001953125,4402.2001953125,3308450000),
    @(45190,4374.35986328125,4375.7001953125,4329.169921875,4330,4330,3662340000),
    @(45191,4341.740234375,4357.39990234375,4316.490234375,4320.06005859375,4320.06005859375,3349570000),
    @(45192,4341.740234375,4357.39990234375,4316.490234375,4320.06005859375,4320.06005859375,3349570000),
    @(45193,4341.740234375,4357.39990234375,4316.490234375,4320.06005859375,4320.06005859375,3349570000),
    @(45194,4310.6201171875,4338.509765625,4302.7001953125,4337.43994140625,4337.43994140625,3195650000),
    @(45195,4312.8798828125,4313.009765625,4265.97998046875,4273.52978515625,4273.52978515625,3472340000),
    @(45196,4282.6298828125,4292.06982421875,4238.6298828125,4274.509765625,4274.509765625,3875880000),
    @(45197,4269.64990234375,4317.27001953125,4264.3798828125,4299.7001953125,4299.7001953125,3846230000),
    @(45198,4328.18017578125,4333.14990234375,4274.85986328125,4288.0498046875,4288.0498046875,3865960000),
    @(45199,4328.18017578125,4333.14990234375,4274.85986328125,4288.0498046875,4288.0498046875,3865960000),
    @(45200,4328.18017578125,4333.14990234375,4274.85986328125,4288.0498046875,4288.0498046875,3865960000),
    @(45201,4284.52001953125,4300.580078125,4260.2099609375,4288.39013671875,4288.39013671875,3938660000),
    @(45202,4269.75,4281.14990234375,4216.4501953125,4229.4501953125,4229.4501953125,3953830000),
    @(45203,4233.830078125,4268.5,4220.47998046875,4263.75,4263.75,3777600000),
    @(45204,4259.31005859375,4267.1298828125,4225.91015625,4258.18994140625,4258.18994140625,3581470000),
    @(45205,4234.7900390625,4324.10009765625,4219.5498046875,4308.5,4308.5,3902030000),
    @(45206,4234.7900390625,4324.10009765625,4219.5498046875,4308.5,4308.5,3902030000),
    @(45207,4234.7900390625,4324.10009765625,4219.5498046875,4308.5,4308.5,3902030000),
    @(45208,4289.02001953125,4341.72998046875,4283.7900390625,4335.66015625,4335.66015625,3174630000),
    @(45209,4339.75,4385.4599609375,4339.64013671875,4358.240234375,4358.240234375,3520240000),
    @(45210,4366.58984375,4378.64013671875,4345.33984375,4376.9501953125,4376.9501953125,3601660000),
    @(45211,4380.93994140625,4385.85009765625,4325.43017578125,4349.60986328125,4349.60986328125,3713140000),
    @(45212,4360.490234375,4377.10009765625,4311.97021484375,4327.77978515625,4327.77978515625,0),
    @(45213,4360.490234375,4377.10009765625,4311.97021484375,4327.77978515625,4327.77978515625,0),
    @(45214,4360.490234375,4377.10009765625,4311.97021484375,4327.77978515625,4327.77978515625,0),
    @(45215,4342.3701171875,4383.330078125,4342.3701171875,4373.6298828125,4373.6298828125,3409960000),
    @(45216,4345.22998046875,4393.56982421875,4337.5400390625,4373.2001953125,4373.2001953125,3794850000),
    @(45217,4357.35009765625,4364.2001953125,4303.83984375,4314.60009765625,4314.60009765625,3686030000),
    @(45218,4321.35986328125,4339.5400390625,4269.68994140625,4278,4278,3969730000),
    @(45219,4273.85009765625,4276.56005859375,4223.02978515625,4224.16015625,4224.16015625,4004030000),
    @(45220,4273.85009765625,4276.56005859375,4223.02978515625,4224.16015625,4224.16015625,4004030000),
    @(45221,4273.85009765625,4276.56005859375,4223.02978515625,4224.16015625,4224.16015625,4004030000),
    @(45222,4210.39990234375,4255.83984375,4189.22021484375,4217.0400390625,4217.0400390625,3776100000)
)

$startRow = 717
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $row.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $row[$j]
    }
}

# Copy the date-column style/format from the last pre-existing data row (A716)
# down through the newly added rows (A717:A754), matching the source workbook's
# per-row date formatting without introducing any new style entries.
$ws.Range("A716").Copy() | Out-Null
$ws.Range("A717:A754").PasteSpecial(-4122) | Out-Null

Write-Host "Appended $($data.Count) rows (717-754) to Sheet1."
